$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Mirriam Makau" row entirely; rows below shift up by one.
$ws.Rows.Item(5).Delete()

# Update the ECO Actual / ECO Balance numbers for the remaining rows
# (kept as text, matching the sheet's existing "N.NN" text convention).
$ws.Range("B2").Value = "'5.00"
$ws.Range("D2").Value = "'5.00"

$ws.Range("B3").Value = "'3.00"
$ws.Range("D3").Value = "'3.00"

$ws.Range("B4").Value = "'1.00"
$ws.Range("D4").Value = "'1.00"

$ws.Range("B6").Value = "'10.00"
$ws.Range("D6").Value = "'10.00"
